# Updated cryptos list values (Price / Volume(1h)) for rows 2-51, columns D and E.
# Source data is text (e.g. "1.000", "28.241.79"), so the Price cells are forced to
# a Text number format while the value is assigned to prevent Excel from silently
# re-interpreting them as numbers/dates, then the style is reset back to Normal so
# no residual cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "28.221.14"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -0.57%  "
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.831.81"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +1.18%  "
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "1.001"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "310.31"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "1.000"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4973"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -3.62%  "
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1009"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +28.10%  "
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3927"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -1.58%  "
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "1.113"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -0.07%  "
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "41.21"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "6.440"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "20.70"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +1.33%  "
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "1.001"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "1.823.03"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "7.343"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001148"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +5.83%  "
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "93.06"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06653"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +1.33%  "
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9998"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +0.01%  "
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "17.27"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "6.017"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "28.270.48"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -0.52%  "
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "11.32"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +1.60%  "
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "2.230"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "158.02"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  -1.79%  "
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "20.82"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "2.040.64"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +0.94%  "
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "2.436"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +1.58%  "
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "126.82"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -1.25%  "
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1052"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -4.17%  "
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "1.039"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -2.78%  "
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "5.607"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "3.600"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -1.81%  "
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06769"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -6.50%  "
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "9.067"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -1.26%  "
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "0.02357"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "0.2154"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -1.48%  "
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "11.46"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -1.44%  "
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "4.990"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -1.29%  "
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6239"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "1.174"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +1.34%  "
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9995"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "13.28"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "0.5947"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -0.86%  "
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "3.689"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -1.14%  "
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "1.272"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -2.93%  "
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "124.14"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -1.25%  "
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "1.954"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +1.15%  "
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "1.182"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -3.27%  "
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "1.127"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +5.05%  "
